$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows holding the 2001-2009 data (rows 2 through 10).
# Using EntireRow.Delete shifts the remaining rows (2010-2018, previously
# rows 11-19) up so they become rows 2-10, matching the target layout.
$ws.Range("A2:E10").EntireRow.Delete()
